# Update the "StructureDefinition-ror-meta-creation-date" workbook for the
# insert release-notes.md commit (f80ed2bb9e1dd81abc71d13817b8a44a756cee80).
#
# Changes:
#   Metadata sheet:
#     Version  : 0.3.0                          -> 0.4.0-snapshot-1
#     Status   : active                         -> draft
#     Date     : 2024-03-13T09:33:00+00:00      -> 2024-05-23T12:16:26+00:00
#     Contact  : No display for ContactDetail   -> ANS (https://esante.gouv.fr)
#
#   Elements sheet:
#     Columns AK (37) and AL (38) - "Mapping: RIM Mapping" and
#     "Mapping: Spécification métier vers l'extension ROR MetaCreationDate" -
#     swap places (header, data and column width).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"
$meta.Range("B6").Value  = "draft"
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# ---------------------------------------------------------------------------
# Elements sheet - swap columns AK (37) and AL (38)
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = 6

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)

    $akValue = $akCell.Value2
    $alValue = $alCell.Value2

    $akCell.Value = $alValue
    $alCell.Value = $akValue
}

# Swap the (best-fit) column widths that go along with the swapped content.
$elements.Columns.Item(37).ColumnWidth = 73.8984375
$elements.Columns.Item(38).ColumnWidth = 24.98046875
